$wb = $excel.ActiveWorkbook

# Update Metadata sheet: Version and Date fields
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("B3").Value = "1.0.0"
$metaSheet.Range("B8").Value = "2025-06-05T14:31:57+02:00"

# Update Elements sheet: clear "Condition(s)" column for Extension.extension and Extension.value[x] rows
$elemSheet = $wb.Worksheets.Item("Elements")
$elemSheet.Range("AI4").Value = ""
$elemSheet.Range("AI6").Value = ""
